$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00653398873036806
$ws.Range("C2").Value = 0.00641409902889342
$ws.Range("D2").Value = 0.00737321664069056
$ws.Range("E2").Value = 0.000059944850737321701893
$ws.Range("F2").Value = 0.999280661791152
$ws.Range("G2").Value = 0.101426687447548
$ws.Range("H2").Value = 0.974523438436638
$ws.Range("I2").Value = 0.000179834552211965
$ws.Range("J2").Value = 0.0787675338688407
$ws.Range("K2").Value = 0.996942812612397
$ws.Range("L2").Value = 0.00035966910442393
$ws.Range("M2").Value = 0.976141949406546
$ws.Range("N2").Value = 0.00167845582064501
$ws.Range("O2").Value = 0.999760220597051
$ws.Range("P2").Value = 0.867761659273468
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.983155496942813
$ws.Range("T2").Value = 0.99910082723894
$ws.Range("U2").Value = 0.000899172761059825
$ws.Range("V2").Value = 0.00749310634216521
$ws.Range("W2").Value = 0.000299724253686608
$ws.Range("X2").Value = 0.00035966910442393

$ws.Range("B3").Value = 0.000059944850737321701893
$ws.Range("C3").Value = 0.000119889701474643
$ws.Range("D3").Value = 0.95971706030452
$ws.Range("E3").Value = 0.999880110298525
$ws.Range("F3").Value = 0.000179834552211965
$ws.Range("G3").Value = 0.000059944850737321701893
$ws.Range("H3").Value = 0.011929025296727
$ws.Range("I3").Value = 0.000179834552211965
$ws.Range("J3").Value = 0.000239779402949287
$ws.Range("K3").Value = 0.00203812492506894
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0.000059944850737321701893
$ws.Range("N3").Value = 0.000059944850737321701893
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.0363265795468169
$ws.Range("Q3").Value = 0.999700275746313
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 0.00185829037285697
$ws.Range("T3").Value = 0.000419613955161252
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0.000059944850737321701893
$ws.Range("W3").Value = 0.000119889701474643
$ws.Range("X3").Value = 0.000059944850737321701893

$ws.Range("B4").Value = 0.993226231866683
$ws.Range("C4").Value = 0.993466011269632
$ws.Range("D4").Value = 0.000659393358110538
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.000419613955161252
$ws.Range("G4").Value = 0.895396235463374
$ws.Range("H4").Value = 0.00929145186428486
$ws.Range("I4").Value = 0.000119889701474643
$ws.Range("J4").Value = 0.918235223594293
$ws.Range("K4").Value = 0.000539503656635895
$ws.Range("L4").Value = 0.99910082723894
$ws.Range("M4").Value = 0.0225392638772329
$ws.Range("N4").Value = 0.998261599328618
$ws.Range("O4").Value = 0.000239779402949287
$ws.Range("P4").Value = 0.083683011629301
$ws.Range("Q4").Value = 0.000059944850737321701893
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.00377652559645126
$ws.Range("T4").Value = 0.000239779402949287
$ws.Range("U4").Value = 0.999040882388203
$ws.Range("V4").Value = 0.991967390001199
$ws.Range("W4").Value = 0.999520441194101
$ws.Range("X4").Value = 0.999460496343364

$ws.Range("D5").Value = 0.0308116532789833
$ws.Range("F5").Value = 0.000119889701474643
$ws.Range("H5").Value = 0.00347680134276466
$ws.Range("I5").Value = 0.999400551492627
$ws.Range("J5").Value = 0.000419613955161252
$ws.Range("K5").Value = 0.00035966910442393
$ws.Range("L5").Value = 0.000299724253686608
$ws.Range("P5").Value = 0.00923150701354754
$ws.Range("Q5").Value = 0.000239779402949287
$ws.Range("S5").Value = 0.0100107900731327
$ws.Range("T5").Value = 0.000119889701474643
$ws.Range("U5").Value = 0.000059944850737321701893
$ws.Range("W5").Value = 0.000059944850737321701893
$ws.Range("X5").Value = 0.000059944850737321701893
